# Be more precise in doc about floating point numbers.
# Replace the answer to question (d) with a more detailed explanation that
# also covers floating point implementation issues.

$d = $word.ActiveDocument

$old = "Es macht keinen Unterschied, da das Minimum genauso funktioniert bei ganzzahligen wie bei nicht-ganzzahligen Zahlen."
$new = "Theoretisch macht es keinen Unterschied, da das Minimum (theoretisch) genauso funktioniert bei ganzzahligen wie bei nicht-ganzzahligen Zahlen. Jedoch ist bei der Implementation von Gleitkomma Zahlen das nicht so der Fall, da die meist dann nicht genau 1.5 sind, sondern 1.50000000004 und würden somit das Ergebnis verfälschen."

$range = $d.Content
$range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
